$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry for 28 марта, continuing the pattern of existing rows.
$ws.Range("B32").Copy()
$ws.Range("B33").PasteSpecial(-4122)

$ws.Range("A33").Value = "28 марта"
$ws.Range("B33").Value = "Интеграция функций cv::Mat в проект и тестирование работы конвейера"

$ws.Range("B34").Select()
